$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.514.50'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '2.382.45'
$ws.Range("E3").Value = '  +6.02%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '233.54'
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("E6").Value = '  +3.39%  '
$ws.Range("D7").Value = '70.03'
$ws.Range("E7").Value = '  +10.93%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.460'
$ws.Range("E9").Value = '  +2.76%  '
$ws.Range("D11").Value = '57.33'
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").Value = '26.34'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.737.08'
$ws.Range("E13").Value = '  +5.97%  '
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("D15").Value = '15.74'
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").Value = '6.24'
$ws.Range("E16").Value = '  +1.95%  '
$ws.Range("D17").Value = '0.853'
$ws.Range("E17").Value = '  +2.89%  '
$ws.Range("D18").Value = '2.373.49'
$ws.Range("E18").Value = '  +5.61%  '
$ws.Range("D19").Value = '43.537.85'
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("D20").Value = '0.0₃0987'
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = '6.35'
$ws.Range("E21").Value = '  +4.72%  '
$ws.Range("D22").Value = '74.15'
$ws.Range("E22").Value = '  +2.10%  '
$ws.Range("D23").Value = '251.18'
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("E24").Value = '  +18.11%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").Value = '2.47'
$ws.Range("E26").Value = '  +2.32%  '
$ws.Range("D27").Value = '2.28'
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("D29").Value = '9.99'
$ws.Range("E29").Value = '  +1.75%  '
$ws.Range("D30").Value = '172.91'
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("D31").Value = '1.55'
$ws.Range("E31").Value = '  +10.09%  '
$ws.Range("E32").Value = '  -8.84%  '
$ws.Range("E33").Value = '  +1.98%  '
$ws.Range("D34").Value = '4.98'
$ws.Range("E34").Value = '  +3.99%  '
$ws.Range("D35").Value = '0.0691'
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("D36").Value = '5.09'
$ws.Range("E36").Value = '  +2.89%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '2.46'
$ws.Range("E37").Value = '  +8.21%  '
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").Value = '6.60'
$ws.Range("E38").Value = '  +3.09%  '
$ws.Range("D39").Value = '3.64'
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("D40").Value = '0.0254'
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '8.94'
$ws.Range("E41").Value = '  +4.52%  '
$ws.Range("B42").Value = 'BinanceUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").Value = '18.56'
$ws.Range("E43").Value = '  +8.89%  '
$ws.Range("D44").Value = '1.20'
$ws.Range("E44").Value = '  +11.89%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '99.36'
$ws.Range("E45").Value = '  +2.03%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = '4.50'
$ws.Range("E46").Value = '  +4.37%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").Value = '1.22'
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("D48").Value = '0.0951'
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("D49").Value = '1.451.09'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("D50").Value = '2.605.93'
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").Value = '2.74'
$ws.Range("E51").Value = '  -0.54%  '